$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 36 and 37 (existing rows 36.. shift down to 38..)
$ws.Range("A36:A37").EntireRow.Insert()

# --- New row 36 ---
$ws.Cells.Item(36, 1).Value = 3
$ws.Cells.Item(36, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(36, 3).Value = "Coquimbo"
$ws.Cells.Item(36, 4).Value = 44879
$ws.Cells.Item(36, 5).Value = 5
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100108
$ws.Cells.Item(36, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(36, 9).Value = 100108004
$ws.Cells.Item(36, 10).Value = "Papaya"
$ws.Cells.Item(36, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 45
$ws.Cells.Item(36, 14).Value = 18000
$ws.Cells.Item(36, 15).Value = 18000
$ws.Cells.Item(36, 16).Value = 18000
$ws.Cells.Item(36, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(36, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(36, 19).Value = 1800
$ws.Cells.Item(36, 20).Value = 10

# --- New row 37 ---
$ws.Cells.Item(37, 1).Value = 3
$ws.Cells.Item(37, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(37, 3).Value = "Coquimbo"
$ws.Cells.Item(37, 4).Value = 44879
$ws.Cells.Item(37, 5).Value = 5
$ws.Cells.Item(37, 6).Value = "Fruta"
$ws.Cells.Item(37, 7).Value = 100108
$ws.Cells.Item(37, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(37, 9).Value = 100108004
$ws.Cells.Item(37, 10).Value = "Papaya"
$ws.Cells.Item(37, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(37, 12).Value = "Segunda"
$ws.Cells.Item(37, 13).Value = 47
$ws.Cells.Item(37, 14).Value = 16000
$ws.Cells.Item(37, 15).Value = 16000
$ws.Cells.Item(37, 16).Value = 16000
$ws.Cells.Item(37, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(37, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(37, 19).Value = 1600
$ws.Cells.Item(37, 20).Value = 10

# Row 56 (previously row 54, shifted down by the insert above) changes Calidad
# from "Primera" to "Segunda".
$ws.Cells.Item(56, 12).Value = "Segunda"

$ws.Range("A1").Select()
